$wb = $excel.ActiveWorkbook
$wsTests = $wb.Worksheets.Item("Tests")
$wsResult = $wb.Worksheets.Item("Result")

# Append the same test cases from "Tests" (rows 2-9, cols A & B) into "Result"
$wsResult.Range("A2").Value = "Framework\InitAllSettings.xaml"
$wsResult.Range("B2").Value = "Success"

$wsResult.Range("A3").Value = "Framework\InitAllApplications.xaml"
$wsResult.Range("B3").Value = "Success"

$wsResult.Range("A4").Value = "Framework\CloseAllApplications.xaml"
$wsResult.Range("B4").Value = "Success"

$wsResult.Range("A5").Value = "Framework\CloseAllApplications.xaml"
$wsResult.Range("B5").Value = "SystemException"

$wsResult.Range("A6").Value = "Framework\InitAllSettings.xaml"
$wsResult.Range("B6").Value = "Success"

$wsResult.Range("A7").Value = "Framework\InitAllSettings.xaml"
$wsResult.Range("B7").Value = "Success"

$wsResult.Range("A8").Value = "Framework\InitAllApplications.xaml"
$wsResult.Range("B8").Value = "Success"

$wsResult.Range("A9").Value = "Framework\CloseAllApplications.xaml"
$wsResult.Range("B9").Value = "Success"

# Update sheet selections to match the latest editing session: the "Tests"
# sheet's cursor moved to B1 (no longer the active tab), while "Result"
# became the active tab with the cursor on B5.
$wsTests.Activate()
$wsTests.Range("B1").Select()

$wsResult.Activate()
$wsResult.Range("B5").Select()
